# Nykaa test case 7: append a third data row (A3 = "600001") to the
# GetAppData test fixture. The cell must hold text (not a parsed number),
# matching how A1/A2 already hold text values ("9123456789", "Kajal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the value as a string-producing formula so Excel stores it with a
# text ("string") cell type instead of inferring a number from the
# digit-only text, then collapse the formula down to its cached text
# value with a Paste-Values so the cell is left as a plain literal (no
# formula, no extra number-format/style) — just like A1 and A2.
$ws.Range("A3").Formula = '="600001"'
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# Excel recomputes $ws.UsedRange / the <dimension> element automatically
# once A3 is populated, extending it from A1:A2 to A1:A3.

# The original workbook ignores the "number stored as text" warning for
# A1:A2; extend that ignored-error range to cover the new row as well, so
# A3 (also digits-as-text) is treated the same way.
$ws.Range("A1:A3").Errors.Item(9).Ignore = $true
